$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (row 6) with a keyword / appID pair, mirroring the
# existing rows: a plain "Normal" style cell in column A and a
# wrap-text styled cell (same formatting as the rest of column B) in
# column B.
$ws.Range("A6").Value = "blockchain free"
$ws.Range("B6").Value = "block.chain.technology"

# Copy the formatting (wrap-text style) from the cell directly above so
# B6 ends up sharing the same cell style as the other column-B cells.
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)

# Move the active selection to the newly added last cell, matching the
# workbook's recorded selection/activeCell.
$ws.Range("B6").Select() | Out-Null
